$p = $ppt.ActivePresentation

# --- 1. Insert a new "Sketches" slide right after "Machine Design" (slide 4) ---
# Duplicate slide 4 so the new slide inherits the exact same layout/placeholder
# naming conventions ("Titel 1" / "Tijdelijke aanduiding voor inhoud 2") used
# throughout the rest of the deck.
$s4 = $p.Slides.Item(4)
$newSlide = $s4.Duplicate().Item(1)

# Retitle it "Sketches".
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Sketches"

# Clear the (duplicated) content placeholder text entirely, leaving a single
# empty paragraph behind (matches a freshly-cleared placeholder).
$newContent = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newContent.Characters(1, 100000).Delete()

# Strip all of the inherited click-animations from the duplicate - the new
# slide has no animation timeline at all.
$newSeq = $newSlide.TimeLine.MainSequence
while ($newSeq.Count -gt 0) {
    $newSeq.Item(1).Delete()
}

# --- 2. Trim the last "Sketch" bullet back out of "Machine Design" (slide 4) ---
$content4 = $s4.Shapes.Item(2).TextFrame.TextRange
$fullLen = $content4.Text.Length
$lastPara = $content4.Paragraphs($content4.Paragraphs().Count)
$content4.Characters($lastPara.Start - 1, $fullLen).Delete()

# Remove the matching (now orphaned) 4th click-animation from Machine Design.
$seq4 = $s4.TimeLine.MainSequence
$seq4.Item($seq4.Count).Delete()
